$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.024.68'
$ws.Range("E2").Value = '  -3.26%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.600.16'
$ws.Range("E3").Value = '  -2.27%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '301.53'
$ws.Range("E6").Value = '  -2.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3783'
$ws.Range("E7").Value = '  -2.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3642'
$ws.Range("E8").Value = '  -4.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '49.76'
$ws.Range("E9").Value = '  -1.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.261'
$ws.Range("E10").Value = '  -4.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.001'
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08133'
$ws.Range("E12").Value = '  -2.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.54'
$ws.Range("E13").Value = '  -4.66%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.595'
$ws.Range("E14").Value = '  -4.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.358'
$ws.Range("E15").Value = '  -5.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001245'
$ws.Range("E16").Value = '  -4.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.605.59'
$ws.Range("E17").Value = '  -1.91%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.98'
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06825'
$ws.Range("E19").Value = '  -1.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.23'
$ws.Range("E20").Value = '  -5.66%  '
$ws.Range("E21").Value = '  -4.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.17'
$ws.Range("E23").Value = '  -2.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.018.45'
$ws.Range("E24").Value = '  -3.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.353'
$ws.Range("E25").Value = '  -3.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.808'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.07'
$ws.Range("E27").Value = '  -3.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.58'
$ws.Range("E28").Value = '  -1.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.230'
$ws.Range("E29").Value = '  -4.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.39'
$ws.Range("E30").Value = '  -1.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.313'
$ws.Range("E31").Value = '  -6.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.840'
$ws.Range("E32").Value = '  -12.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.784.90'
$ws.Range("E33").Value = '  -1.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9621'
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07589'
$ws.Range("E35").Value = '  -4.45%  '
$ws.Range("E36").Value = '  -0.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.263'
$ws.Range("E37").Value = '  -4.18%  '
$ws.Range("E38").Value = '  -5.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2533'
$ws.Range("E39").Value = '  -4.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08895'
$ws.Range("E40").Value = '  -1.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.369'
$ws.Range("E41").Value = '  -3.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7026'
$ws.Range("E42").Value = '  -5.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.43'
$ws.Range("E43").Value = '  -5.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.16'
$ws.Range("E44").Value = '  -8.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6631'
$ws.Range("E45").Value = '  -3.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.305'
$ws.Range("E47").Value = '  -3.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.992'
$ws.Range("E48").Value = '  -1.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '131.56'
$ws.Range("E49").Value = '  -1.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07909'
$ws.Range("E50").Value = '  -3.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.230'
$ws.Range("E51").Value = '  +1.36%  '
